$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.502.45"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "1.571.33"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("E6").Value = "  -3.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.26%  "
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0587"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").Value = "1.789.97"
$ws.Range("E12").Value = "  -3.97%  "
$ws.Range("D13").Value = "1.555.48"
$ws.Range("E13").Value = "  -4.74%  "
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("E15").Value = "  -6.04%  "
$ws.Range("D16").Value = "27.456.88"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.35%  "
$ws.Range("D20").Value = "0.0₃0686"
$ws.Range("E20").Value = "  -4.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.16%  "
$ws.Range("E24").Value = "  -4.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("E29").Value = "  -4.67%  "
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("E31").Value = "  -3.96%  "
$ws.Range("E32").Value = "  -5.96%  "
$ws.Range("D33").Value = "1.356.78"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("E34").Value = "  -5.73%  "
$ws.Range("E35").Value = "  -5.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.971"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("E39").Value = "  -4.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.807"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.971"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").Value = "1.704.80"
$ws.Range("E47").Value = "  -3.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0964"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.12%  "
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").Value = "0.0₇0964"
$ws.Range("E51").Value = "  -6.86%  "
